# ---------------------------------------------------------------------------
# Append a new data row (row 13) for symbol AAME to the stock-overview sheet.
#
# Commit message: "Changed the deleteNoneUpdateSymbols func in order to drop
# the ExcelWriter Object" - the upstream generator no longer runs its values
# through pandas/openpyxl's ExcelWriter (which used to infer numeric dtypes),
# so every field in the appended row - except the fixed "downloaded at" audit
# columns - is now written out as a literal text value, regardless of whether
# it looks like a number ("0", "8177", "-0.05", ...).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRow = 13

# Temporarily force text storage for columns A:AY so values such as "0" or
# "8177" are not auto-coerced into numbers by Excel's input parser.
$textRange = "A" + $dataRow + ":AY" + $dataRow
$ws.Range($textRange).NumberFormat = "@"

$ws.Cells.Item($dataRow, 1).Value = "AAME"
$ws.Cells.Item($dataRow, 2).Value = "Common Stock"
$ws.Cells.Item($dataRow, 3).Value = "Atlantic American Corporation"
$ws.Cells.Item($dataRow, 4).Value = "Atlantic American Corporation provides life and health insurance and property and casualty products in the United States. The company is headquartered in Atlanta, Georgia."
$ws.Cells.Item($dataRow, 5).Value = "8177"
$ws.Cells.Item($dataRow, 6).Value = "NASDAQ"
$ws.Cells.Item($dataRow, 7).Value = "USD"
$ws.Cells.Item($dataRow, 8).Value = "USA"
$ws.Cells.Item($dataRow, 9).Value = "FINANCE"
$ws.Cells.Item($dataRow, 10).Value = "LIFE INSURANCE"
$ws.Cells.Item($dataRow, 11).Value = "4370 PEACHTREE RD NE, ATLANTA, GA, US"
$ws.Cells.Item($dataRow, 12).Value = "December"
$ws.Cells.Item($dataRow, 13).Value = "2024-03-31"
$ws.Cells.Item($dataRow, 14).Value = "32146000"
$ws.Cells.Item($dataRow, 15).Value = "2698000"
$ws.Cells.Item($dataRow, 16).Value = "None"
$ws.Cells.Item($dataRow, 17).Value = "0"
$ws.Cells.Item($dataRow, 18).Value = "5.04"
$ws.Cells.Item($dataRow, 19).Value = "0.02"
$ws.Cells.Item($dataRow, 20).Value = "0.0127"
$ws.Cells.Item($dataRow, 21).Value = "-0.05"
$ws.Cells.Item($dataRow, 22).Value = "9.19"
$ws.Cells.Item($dataRow, 23).Value = "-0.0039"
$ws.Cells.Item($dataRow, 24).Value = "-0.0351"
$ws.Cells.Item($dataRow, 25).Value = "0.0037"
$ws.Cells.Item($dataRow, 26).Value = "-0.007"
$ws.Cells.Item($dataRow, 27).Value = "187521000"
$ws.Cells.Item($dataRow, 28).Value = "20364000"
$ws.Cells.Item($dataRow, 29).Value = "-0.05"
$ws.Cells.Item($dataRow, 30).Value = "-0.784"
$ws.Cells.Item($dataRow, 31).Value = "0.016"
$ws.Cells.Item($dataRow, 32).Value = "None"
$ws.Cells.Item($dataRow, 33).Value = "0"
$ws.Cells.Item($dataRow, 34).Value = "0"
$ws.Cells.Item($dataRow, 35).Value = "0"
$ws.Cells.Item($dataRow, 36).Value = "0"
$ws.Cells.Item($dataRow, 37).Value = "0"
$ws.Cells.Item($dataRow, 38).Value = "-"
$ws.Cells.Item($dataRow, 39).Value = "-"
$ws.Cells.Item($dataRow, 40).Value = "0.239"
$ws.Cells.Item($dataRow, 41).Value = "0.417"
$ws.Cells.Item($dataRow, 42).Value = "0.285"
$ws.Cells.Item($dataRow, 43).Value = "-"
$ws.Cells.Item($dataRow, 44).Value = "0.426"
$ws.Cells.Item($dataRow, 45).Value = "2.974"
$ws.Cells.Item($dataRow, 46).Value = "1.477"
$ws.Cells.Item($dataRow, 47).Value = "1.695"
$ws.Cells.Item($dataRow, 48).Value = "2.088"
$ws.Cells.Item($dataRow, 49).Value = "20399800"
$ws.Cells.Item($dataRow, 50).Value = "2024-04-26"
$ws.Cells.Item($dataRow, 51).Value = "2024-04-11"

# Restore the plain "Normal" style on the text columns so no lingering
# text-number-format style is left behind on the cells.
$ws.Range($textRange).Style = "Normal"

# Download audit columns stay numeric: timestamp (AZ, same style as the rows
# above it) and quarter number (BA).
$azCell = "AZ" + $dataRow
$ws.Range("AZ12").Copy()
$ws.Range($azCell).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item($dataRow, 52).Value = 45498.85622998936
$ws.Cells.Item($dataRow, 53).Value = 3

